$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------
# 1. Insert a brand new row above the old row 5 ("dioda przetwornicy").
#    This pushes old rows 5..15 down to 6..16 and leaves row 4 (the
#    old "cewka 47uH 10A" row) in place so it can be overwritten with
#    new content below.
# ---------------------------------------------------------------
$ws.Range("A5").EntireRow.Insert()

# ---------------------------------------------------------------
# 2. Add the new comment in H2
# ---------------------------------------------------------------
$ws.Range("H2").Value2 = "zamówiona partia z chin, najwyzej wylut z przetwornicy z allegro"

# ---------------------------------------------------------------
# 3. Rebuild row 4 completely with the new "cewka 47uH 6.8A SMD" part
# ---------------------------------------------------------------
$ws.Range("A4:J4").ClearContents()
$ws.Range("B4").Value2 = "cewka 47uH 6.8A SMD"
$ws.Range("C4").Value2 = "ETQP5M470YFC"
$ws.Range("D4").Value2 = 2
$ws.Range("F4").Value2 = "ETQP5M470YFC"
$ws.Range("G4").Value2 = 3.11
$ws.Range("H4").Value2 = "imo lepiej smd"
$ws.Range("I4").Value2 = "Inductor_SMD:L_Bourns-SRN8040_8x8.15mm"

# ---------------------------------------------------------------
# 4. Fill in the new row 5 with the "cewka 47uH 5A THD" part
# ---------------------------------------------------------------
$ws.Range("B5").Value2 = "cewka 47uH 5A THD"
$ws.Range("C5").Value2 = "DPO-5.0-47"
$ws.Range("D5").Value2 = 2
$ws.Range("F5").Value2 = "DPU047A5"
$ws.Range("G5").Value2 = 6
$ws.Range("H5").Value2 = "trzeba stwierdzić czy 5A wystarczy (według noty przetwornicy tak)"
$ws.Range("I5").Value2 = "Inductor_THT:L_Toroid_Vertical_L25.4mm_W14.7mm_P12.20mm_Vishay_TJ5"

# ---------------------------------------------------------------
# 5. E3:E5 share the "=D*4" formula (mirrors the rest of column E)
# ---------------------------------------------------------------
$ws.Range("E3:E5").Formula = "=D3*4"

# ---------------------------------------------------------------
# 6. Style rows 4 and 5: yellow fill for every used cell, and the
#    currency number format (kept) + yellow fill for column G.
# ---------------------------------------------------------------
$rowsToFill = @("B4","C4","D4","E4","F4","H4","I4","B5","C5","D5","E5","F5","H5","I5")
foreach ($addr in $rowsToFill) {
    $ws.Range($addr).Interior.Color = 65535
}
$ws.Range("G4").NumberFormat = "#,##0.00\ ""zł"""
$ws.Range("G4").Interior.Color = 65535
$ws.Range("G5").NumberFormat = "#,##0.00\ ""zł"""
$ws.Range("G5").Interior.Color = 65535

# ---------------------------------------------------------------
# 7. The B540C-13-F "?" note becomes "? / 3220"
# ---------------------------------------------------------------
$ws.Range("I6").Value2 = "? / 3220"

# ---------------------------------------------------------------
# 8. The "0.14" value in G7 (dioda sterowania price) becomes numeric
# ---------------------------------------------------------------
$ws.Range("G7").Value2 = 0.14

# ---------------------------------------------------------------
# 9. Fix hyperlinks: only J2, J3 and J6 (shifted from old J5) remain
# ---------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "https://www.elecrow.com/download/XL4015_datasheet.pdf")
$ws.Hyperlinks.Add($ws.Range("J3"), "https://www.tme.eu/Document/c98656c49a036767b89c9bb93e3dda4a/mcp6001_2_4.pdf")
$ws.Hyperlinks.Add($ws.Range("J6"), "https://www.tme.eu/pl/details/b540c-13-f/diody-schottky-smd/diodes-incorporated/")

# ---------------------------------------------------------------
# 10. Column / view cosmetics
# ---------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(8).ColumnWidth = 59.28515625
$ws.Range("H28").Select()

Write-Host "done"
